$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 485, pushing existing rows 485:508 down to 486:509.
$ws.Rows.Item(485).Insert()

# Populate the newly inserted row 485 with the new data record.
$ws.Range("A485").Value = 3
$ws.Range("B485").Value = "Femacal de La Calera"
$ws.Range("C485").Value = "Coquimbo"
$ws.Range("D485").Value = 44753
$ws.Range("D485").NumberFormat = $ws.Range("D486").NumberFormat
$ws.Range("E485").Value = 5
$ws.Range("F485").Value = 100112037
$ws.Range("G485").Value = "Cebollín"
$ws.Range("H485").Value = "Sin especificar"
$ws.Range("I485").Value = "Primera"
$ws.Range("J485").Value = 190
$ws.Range("K485").Value = 6500
$ws.Range("L485").Value = 7000
$ws.Range("M485").Value = 6684
$ws.Range("N485").Value = "$/paquete 36 unidades"
$ws.Range("O485").Value = "Provincia de Quillota"
$ws.Range("P485").Value = 186
$ws.Range("Q485").Value = 36
$ws.Range("R485").Value = "Hortaliza"
